# Applies the "updated for test step execution" edit:
#  - Remove the DASHBOARD and LOGOUT worksheets, keeping only LOGIN
#  - Insert a new row (row 2) into the LOGIN sheet with a LOGIN/LOGIN/LOGIN step
#  - Update the dimension/selection accordingly

$wb = $excel.ActiveWorkbook

# Turn off alerts so sheet deletion doesn't prompt
$excel.DisplayAlerts = $false

# Remove the DASHBOARD and LOGOUT sheets, leaving only LOGIN
$dashboard = $wb.Worksheets.Item("DASHBOARD")
$dashboard.Delete()

$logout = $wb.Worksheets.Item("LOGOUT")
$logout.Delete()

$excel.DisplayAlerts = $true

# Work on the remaining LOGIN sheet
$ws = $wb.Worksheets.Item("LOGIN")

# Insert a new row above the current row 2 (shifts the SEND/SEND/CLICK rows down)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the LOGIN step
$ws.Range("A2").Value = "LOGIN"
$ws.Range("B2").Value = "LOGIN"
$ws.Range("E2").Value = "LOGIN"

# Update the selected cell to reflect the new active cell E3
$ws.Range("E3").Select()
